$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (shifts existing rows 9-13 down to 10-14)
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the latest weekly record
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44669
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112001
$ws.Range("G9").Value = "Berenjena"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 130
$ws.Range("K9").Value = 4500
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = 4750
$ws.Range("N9").Value = "$/caja 60 unidades"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 79
$ws.Range("Q9").Value = 60
$ws.Range("R9").Value = "Hortaliza"
